$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H. This shifts the old "Ket:" note block
# (H5:I6) one column to the right (H->I, I->J), matching the target
# layout where two new header cells ("Kwitansi" / "Tipe Formulir") are
# inserted ahead of the existing note column.
$ws.Columns("H").Insert()

# G5 / H5: new header cells, bordered box + vertically centered text
# (no fill, no bold, no wrap/h-center) - build that style by copying an
# existing box-bordered, vertically centered cell (E5) and stripping its
# bold font.
$ws.Range("E5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Font.Bold = $False
$ws.Range("G5").Value = "Kwitansi"

$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "Tipe Formulir"

# New row-6 cells under the new columns get the same plain box-border
# style already used by B6/C6 in that row.
$ws.Range("B6").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)

# Column width tweaks: F widens (fixed width instead of autosize) to fit
# the new "Kwitansi" header, and the new column H gets an explicit width
# for "Tipe Formulir".
$ws.Columns("F").ColumnWidth = 22.1
$ws.Columns("H").ColumnWidth = 17.5

# Selection left where the author's resave recorded it.
$ws.Range("F10").Select()
